$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the touched cells to remain text, matching the workbook's
# existing convention of storing all values (including numeric-looking
# ones) as strings.
$ws.Range("C4:C8").NumberFormat = "@"
$ws.Range("E4:E8").NumberFormat = "@"

# Row 4
$ws.Range("C4").Value = "1"
$ws.Range("E4").Value = "9.25"

# Row 5
$ws.Range("C5").Value = "1"
$ws.Range("E5").Value = "23.75"

# Row 6
$ws.Range("C6").Value = "1"
$ws.Range("E6").Value = "13.00"

# Row 7
$ws.Range("C7").Value = "1"
$ws.Range("E7").Value = "13.00"

# Row 8
$ws.Range("C8").Value = "1"
$ws.Range("E8").Value = "9.25"
